$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The worksheet is protected (no password needed to unprotect in this
# environment, but we pass the known password to be safe and to allow
# re-protecting with it afterwards).
$ws.Unprotect("D382")

# 1) Update the "as of" date embedded in the confidential disclaimer text (A7).
$disclaimer = $ws.Range("A7").Value2
$ws.Range("A7").Value2 = $disclaimer.Replace("2021-05-11", "2021-05-12")

# 2) Update the weight / percent-change figures for rows 2-4.
$ws.Range("D2").Value2 = 0.8467733283745413
$ws.Range("E2").Value2 = -0.01596351197263401

$ws.Range("D3").Value2 = 0.1532266716254586
$ws.Range("E3").Value2 = -0.0276325616131442

$ws.Range("E4").Value2 = -0.01775152161008176

# Restore sheet protection as it was.
$ws.Protect("D382")
